$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.469.97"
Set-TextValue $ws.Range("E2") "  +2.38%  "

Set-TextValue $ws.Range("D3") "1.826.70"
Set-TextValue $ws.Range("E3") "  +1.59%  "

Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.02%  "

Set-TextValue $ws.Range("D5") "315.13"
Set-TextValue $ws.Range("E5") "  -0.20%  "

Set-TextValue $ws.Range("E6") "  -0.12%  "

Set-TextValue $ws.Range("D7") "0.5114"
Set-TextValue $ws.Range("E7") "  -4.74%  "

Set-TextValue $ws.Range("D8") "0.3926"
Set-TextValue $ws.Range("E8") "  +1.75%  "

Set-TextValue $ws.Range("D9") "0.07728"
Set-TextValue $ws.Range("E9") "  +3.90%  "

Set-TextValue $ws.Range("E10") "  +1.15%  "

Set-TextValue $ws.Range("D11") "1.112"
Set-TextValue $ws.Range("E11") "  +2.31%  "

Set-TextValue $ws.Range("D12") "21.06"
Set-TextValue $ws.Range("E12") "  +3.57%  "

Set-TextValue $ws.Range("D13") "6.278"
Set-TextValue $ws.Range("E13") "  +1.15%  "

Set-TextValue $ws.Range("D14") "7.581"
Set-TextValue $ws.Range("E14") "  +0.90%  "

Set-TextValue $ws.Range("D15") "1.002"
Set-TextValue $ws.Range("E15") "  +0.04%  "

Set-TextValue $ws.Range("D16") "1.826.71"
Set-TextValue $ws.Range("E16") "  +2.19%  "

Set-TextValue $ws.Range("D17") "93.48"
Set-TextValue $ws.Range("E17") "  +5.66%  "

Set-TextValue $ws.Range("D18") "0.00001083"
Set-TextValue $ws.Range("E18") "  +2.20%  "

Set-TextValue $ws.Range("D19") "0.06628"
Set-TextValue $ws.Range("E19") "  +2.07%  "

Set-TextValue $ws.Range("D20") "17.71"
Set-TextValue $ws.Range("E20") "  +2.37%  "

Set-TextValue $ws.Range("E21") "  -0.12%  "

Set-TextValue $ws.Range("E22") "  +2.41%  "

Set-TextValue $ws.Range("D23") "28.478.05"
Set-TextValue $ws.Range("E23") "  +2.31%  "

Set-TextValue $ws.Range("D24") "11.11"
Set-TextValue $ws.Range("E24") "  -0.09%  "

Set-TextValue $ws.Range("D25") "2.259"
Set-TextValue $ws.Range("E25") "  +7.81%  "

Set-TextValue $ws.Range("D26") "157.09"
Set-TextValue $ws.Range("E26") "  +0.69%  "

Set-TextValue $ws.Range("D27") "2.441"
Set-TextValue $ws.Range("E27") "  +4.51%  "

Set-TextValue $ws.Range("D28") "20.61"
Set-TextValue $ws.Range("E28") "  +1.52%  "

Set-TextValue $ws.Range("D29") "2.035.86"
Set-TextValue $ws.Range("E29") "  +2.00%  "

Set-TextValue $ws.Range("D30") "124.78"
Set-TextValue $ws.Range("E30") "  +2.58%  "

Set-TextValue $ws.Range("D31") "1.131"
Set-TextValue $ws.Range("E31") "  +1.05%  "

Set-TextValue $ws.Range("D32") "0.1098"
Set-TextValue $ws.Range("E32") "  +0.42%  "

Set-TextValue $ws.Range("D33") "5.647"
Set-TextValue $ws.Range("E33") "  +2.45%  "

Set-TextValue $ws.Range("D34") "3.677"
Set-TextValue $ws.Range("E34") "  +0.64%  "

Set-TextValue $ws.Range("D35") "0.07147"
Set-TextValue $ws.Range("E35") "  +1.27%  "

Set-TextValue $ws.Range("D36") "0.2234"
Set-TextValue $ws.Range("E36") "  +1.63%  "

Set-TextValue $ws.Range("D37") "8.973"
Set-TextValue $ws.Range("E37") "  +5.86%  "

Set-TextValue $ws.Range("D38") "0.02322"
Set-TextValue $ws.Range("E38") "  +1.96%  "

Set-TextValue $ws.Range("D39") "5.145"
Set-TextValue $ws.Range("E39") "  +1.74%  "

Set-TextValue $ws.Range("D40") "0.6236"
Set-TextValue $ws.Range("E40") "  +2.12%  "

Set-TextValue $ws.Range("D41") "11.26"
Set-TextValue $ws.Range("E41") "  -0.44%  "

Set-TextValue $ws.Range("D42") "1.188"
Set-TextValue $ws.Range("E42") "  +2.21%  "

Set-TextValue $ws.Range("E43") "  -0.17%  "

Set-TextValue $ws.Range("E44") "  -1.08%  "

Set-TextValue $ws.Range("D45") "13.45"
Set-TextValue $ws.Range("E45") "  +1.06%  "

Set-TextValue $ws.Range("D46") "0.5890"
Set-TextValue $ws.Range("E46") "  +3.31%  "

Set-TextValue $ws.Range("D47") "3.705"
Set-TextValue $ws.Range("E47") "  +0.74%  "

Set-TextValue $ws.Range("D48") "124.29"
Set-TextValue $ws.Range("E48") "  -0.45%  "

Set-TextValue $ws.Range("D49") "1.977"
Set-TextValue $ws.Range("E49") "  +3.54%  "

Set-TextValue $ws.Range("D50") "1.182"
Set-TextValue $ws.Range("E50") "  +0.78%  "

Set-TextValue $ws.Range("D51") "0.06923"
Set-TextValue $ws.Range("E51") "  +1.96%  "
